# Fixed typo in session 3 slides.
#
# 1) Slide 1 ("Rectangle 3" textbox): the date was split across two runs
#    ("Thursday, September " + "19, 2013") -- merge them into a single run.
# 2) Slide 22 ("Content Placeholder 2"): "number between 0 and 3: " had a
#    typo (should read "... 0 and 2: ") -- fix the digit and, in doing so,
#    the run gets split right after "and ".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Edit 1: slide 1 - merge the "Thursday, September " / "19, 2013" runs.
# ---------------------------------------------------------------------
$s1  = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange

# Find the paragraph that holds the date (works even if paragraph order
# ever shifts, since we search by content rather than a hard-coded index).
$dateParaIndex = -1
for ($i = 1; $i -le $tr1.Paragraphs(1, -1).Count; $i++) {
    $candidate = $tr1.Paragraphs($i, 1)
    if ($candidate.Text -like "Thursday, September*2013*") {
        $dateParaIndex = $i
        break
    }
}
if ($dateParaIndex -eq -1) { $dateParaIndex = 3 }

$datePara = $tr1.Paragraphs($dateParaIndex, 1)
$dateText = $datePara.Text
$splitAt = $dateText.IndexOf("19, 2013")

$firstRun  = $tr1.Characters($datePara.Start, $splitAt)
$secondRun = $tr1.Characters($datePara.Start + $splitAt, $dateText.Length - $splitAt)

$secondText = $secondRun.Text
$secondRun.Delete() | Out-Null
$firstRun.InsertAfter($secondText) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: slide 22 - "number between 0 and 3: " -> "number between 0 and 2: "
# ---------------------------------------------------------------------
$s22  = $p.Slides.Item(22)
$sh22 = $s22.Shapes.Item(2)
$tr22 = $sh22.TextFrame.TextRange

$exerciseParaIndex = -1
for ($i = 1; $i -le $tr22.Paragraphs(1, -1).Count; $i++) {
    $candidate = $tr22.Paragraphs($i, 1)
    if ($candidate.Text -like "*number between 0 and*") {
        $exerciseParaIndex = $i
        break
    }
}
if ($exerciseParaIndex -eq -1) { $exerciseParaIndex = 7 }

$exercisePara = $tr22.Paragraphs($exerciseParaIndex, 1)
$exerciseText = $exercisePara.Text
$relIdx = $exerciseText.IndexOf("3: ")
$absStart = $exercisePara.Start + $relIdx

# Replace the whole "3: " token (not just the digit) so the run splits
# cleanly into "number between 0 and " + "2: " (matching how PowerPoint
# re-segments runs on a targeted retype of that token).
$oldRange = $tr22.Characters($absStart, 3)
$oldRange.Text = "2: "

Write-Output "Applied session 3 typo fixes."
